$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column J, row 3
$ws.Range("J3").Value = "After cleanup"

# New data rows for column J (mirrors column I structure)
$ws.Range("J5").Value = $ws.Range("I5").Value2
$ws.Range("J6").Value = $ws.Range("I6").Value2
$ws.Range("J7").Value = "1,20,511"

$ws.Range("J8").Value = $ws.Range("I8").Value2
$ws.Range("J9").Value = $ws.Range("I9").Value2
$ws.Range("J10").Value = "20,2,966"

$ws.Range("J11").Value = $ws.Range("I11").Value2
$ws.Range("J12").Value = $ws.Range("I12").Value2
$ws.Range("J13").Value = "1,20,500"

$ws.Range("J14").Value = $ws.Range("I14").Value2
$ws.Range("J15").Value = $ws.Range("I15").Value2
$ws.Range("J16").Value = "20,2,958"

$ws.Range("J17").Value = $ws.Range("I17").Value2
$ws.Range("J18").Value = $ws.Range("I18").Value2
$ws.Range("J19").Value = "1,20,495"

$ws.Range("J20").Value = $ws.Range("I20").Value2
$ws.Range("J21").Value = $ws.Range("I21").Value2
$ws.Range("J22").Value = "20,2,955"

# Update selection to reflect new active cell
$ws.Range("K7").Select()
